# Apply updated bus voltage magnitude (vm_pu) results for the 380 kV case
# (commit: "case with 380 kV done") - bus 0 setpoint lowered from 1.05 to 1.02 pu,
# and the rest of the load-flow solution columns (B-F, I-N) are updated accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
  2 = @{ "B"=1.02; "C"=1.046458953949363; "D"=1.05085567420715; "E"=1.050126891390883; "F"=1.058919810969283; "I"=1.043013168987303; "J"=1.051513238183965; "K"=1.053608792000661; "L"=1.052882033625601; "M"=1.061650729357352; "N"=1.02095638341446 }
  3 = @{ "B"=1.02; "C"=1.04754372611421; "D"=1.05170569495475; "E"=1.051161117473204; "F"=1.059964188561709; "I"=1.043311086746218; "J"=1.052245420299244; "K"=1.054271246094916; "L"=1.053728071762166; "M"=1.062508651516115; "N"=1.021206120846143 }
  4 = @{ "B"=1.02; "C"=1.048245792057101; "D"=1.052255839449221; "E"=1.051830813601488; "F"=1.060640522624598; "I"=1.043502741708715; "J"=1.052718760902137; "K"=1.054699378149174; "L"=1.054275392167799; "M"=1.063063730501904; "N"=1.021367420852036 }
  5 = @{ "B"=1.02; "C"=1.04854097581055; "D"=1.052487149522707; "E"=1.052112469255426; "F"=1.060924985343283; "I"=1.043583045753587; "J"=1.052917650204392; "K"=1.05487924016709; "L"=1.054505456136956; "M"=1.063297072565688; "N"=1.021435160169379 }
  6 = @{ "B"=1.02; "C"=1.048590540537086; "D"=1.052525989209682; "E"=1.052159767241783; "F"=1.060972755611577; "I"=1.043596513454787; "J"=1.052951038545463; "K"=1.054909432473266; "L"=1.054544083154641; "M"=1.06333625094827; "N"=1.021446529719035 }
  7 = @{ "B"=1.02; "C"=1.048249736178009; "D"=1.052258930111019; "E"=1.051834576644858; "F"=1.060644323111483; "I"=1.043503815787307; "J"=1.05272141887591; "K"=1.054701781964481; "L"=1.05427846640999; "M"=1.063066848481875; "N"=1.021368326268081 }
  8 = @{ "B"=1.02; "C"=1.046825527782304; "D"=1.051142916727247; "E"=1.050476312897503; "F"=1.059272648834019; "I"=1.043114083059448; "J"=1.051760771570202; "K"=1.053832779021877; "L"=1.053167981537207; "M"=1.061940679226216; "N"=1.021040844716883 }
  9 = @{ "B"=1.02; "C"=1.044316992664666; "D"=1.049177327923603; "E"=1.048086582033449; "F"=1.056859821393666; "I"=1.042418768297724; "J"=1.050064697647502; "K"=1.052297510377319; "L"=1.05121022945634; "M"=1.059955820302681; "N"=1.020461513549286 }
  10 = @{ "B"=1.02; "C"=1.042645353529063; "D"=1.047867609129508; "E"=1.046495924115091; "F"=1.055254130996995; "I"=1.041949477839475; "J"=1.048931772492606; "K"=1.051271335187909; "L"=1.049904436420728; "M"=1.058632318482804; "N"=1.020073776542618 }
  11 = @{ "B"=1.02; "C"=1.041921679468916; "D"=1.047300650173607; "E"=1.045807744062996; "F"=1.054559529620408; "I"=1.041744907421843; "J"=1.048440679331346; "K"=1.050826359757834; "L"=1.04933886546459; "M"=1.058059166562746; "N"=1.019905523917659 }
  12 = @{ "B"=1.02; "C"=1.041652897617227; "D"=1.047090080333078; "E"=1.045552211212508; "F"=1.05430162490831; "I"=1.041668715847797; "J"=1.048258185872788; "K"=1.050660980718848; "L"=1.049128764004664; "M"=1.057846262372897; "N"=1.019842973404163 }
  13 = @{ "B"=1.02; "C"=1.041710551209888; "D"=1.047135247190706; "E"=1.045607019901627; "F"=1.054356941766026; "I"=1.041685068476354; "J"=1.048297334942493; "K"=1.050696459398636; "L"=1.049173832533499; "M"=1.057891931498922; "N"=1.01985639314384 }
  14 = @{ "B"=1.02; "C"=1.041899461413354; "D"=1.047283243914784; "E"=1.045786619845933; "F"=1.0545382090725; "I"=1.04173861358477; "J"=1.04842559598938; "K"=1.050812691419905; "L"=1.049321498886688; "M"=1.058041568040928; "N"=1.019900354573763 }
  15 = @{ "B"=1.02; "C"=1.042015858242365; "D"=1.04737443276285; "E"=1.045897288939504; "F"=1.054649907249087; "I"=1.041771577327063; "J"=1.048504611303709; "K"=1.050884293169949; "L"=1.049412477920189; "M"=1.058133762721883; "N"=1.0199274335089 }
  16 = @{ "B"=1.02; "C"=1.042693384641122; "D"=1.047905239675366; "E"=1.04654160868423; "F"=1.055300243588122; "I"=1.041963025754427; "J"=1.048964353572186; "K"=1.051300853368851; "L"=1.049941968250638; "M"=1.058670355278903; "N"=1.02008493533848 }
  17 = @{ "B"=1.02; "C"=1.043118420563494; "D"=1.048238243273496; "E"=1.046945930090972; "F"=1.055708362638165; "I"=1.042082750962189; "J"=1.049252596025139; "K"=1.051561980635749; "L"=1.050274062260603; "M"=1.059006927704941; "N"=1.020183635764904 }
  18 = @{ "B"=1.02; "C"=1.043366351988866; "D"=1.048432493860141; "E"=1.047181820407196; "F"=1.055946476523342; "I"=1.04215245299092; "J"=1.049420671928243; "K"=1.05171423050208; "L"=1.050467752229902; "M"=1.059203238251671; "N"=1.020241171319177 }
  19 = @{ "B"=1.02; "C"=1.043450892788874; "D"=1.048498730865093; "E"=1.047262262491441; "F"=1.056027678220162; "I"=1.042176197249891; "J"=1.049477972813084; "K"=1.051766133382267; "L"=1.050533792988443; "M"=1.059270173980015; "N"=1.020260783562304 }
  20 = @{ "B"=1.02; "C"=1.043072816667551; "D"=1.048202513562692; "E"=1.046902544388905; "F"=1.055664568608269; "I"=1.042069919198591; "J"=1.049221675618098; "K"=1.051533970486494; "L"=1.050238433252632; "M"=1.058970817336463; "N"=1.02017304973539 }
  21 = @{ "B"=1.019999999999999; "C"=1.041843831439339; "D"=1.047239661904991; "E"=1.045733729713924; "F"=1.054484827575766; "I"=1.041722851544285; "J"=1.04838782850585; "K"=1.050778466618655; "L"=1.04927801546803; "M"=1.057997504086213; "N"=1.019887410518302 }
  22 = @{ "B"=1.02; "C"=1.0410712519237; "D"=1.046634417590432; "E"=1.044999357844251; "F"=1.053743662898705; "I"=1.041503450195867; "J"=1.047863095183958; "K"=1.050302899686399; "L"=1.048674028271148; "M"=1.057385485288453; "N"=1.019707505434342 }
  23 = @{ "B"=1.02; "C"=1.041480798613484; "D"=1.046955255824557; "E"=1.045388614110538; "F"=1.054136512823527; "I"=1.041619871434909; "J"=1.048141309934436; "K"=1.050555058999639; "L"=1.048994226084402; "M"=1.057709933541065; "N"=1.019802906105732 }
  24 = @{ "B"=1.02; "C"=1.04309342305803; "D"=1.04821865823393; "E"=1.046922148345497; "F"=1.055684357044485; "I"=1.042075717726967; "J"=1.049235647376673; "K"=1.051546627256414; "L"=1.050254532513445; "M"=1.058987134077453; "N"=1.020177833213755 }
  25 = @{ "B"=1.02; "C"=1.044965380940761; "D"=1.049685362601717; "E"=1.048703945254007; "F"=1.057483091686957; "I"=1.042599537216838; "J"=1.050503563039622; "K"=1.052694884144613; "L"=1.051716465767883; "M"=1.060469000894461; "N"=1.020611552362659 }
}

foreach ($rowNum in $newValues.Keys) {
  $rowData = $newValues[$rowNum]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$rowNum").Value = $rowData[$col]
  }
}
